$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new hyperlink in B8 (text + target both https://www.hostpapa.com),
# mirroring the existing hyperlink already in B7.
$ws.Hyperlinks.Add($ws.Range("B8"), "https://www.hostpapa.com") | Out-Null

# Apply the built-in Hyperlink cell style, same as the existing link in B7.
$ws.Range("B8").Style = "Hyperlink"

# The recorded test run ends with the user having clicked/selected G9
# (e.g. after following the link in the browser and returning to Excel).
$ws.Range("G9").Select() | Out-Null
